# Refresh the "cryptos" price/volume snapshot (GitHub Actions data pull).
# For each affected row we rewrite the changed columns only:
#   D = Price (kept as literal text, e.g. "60.978.71", never coerced to a number)
#   E = Volume(1h) (kept as the "  +x.xx%  " padded text)
#   B/C = Coin name / link (rows 39-40: VeChain and Hedera swapped rank order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='60.978.71'; E='  +4.26%  ' }  # row 2: Bitcoin
    @{ Row=3; D='3.249.11'; E='  +3.08%  ' }  # row 3: Ethereum
    @{ Row=4; D='1.00'; E='  +0.02%  ' }  # row 4: TetherUSD
    @{ Row=5; D='545.07'; E='  +3.12%  ' }  # row 5: BNB
    @{ Row=6; D='146.73'; E='  +5.21%  ' }  # row 6: Solana
    @{ Row=7; E='  +0.09%  ' }  # row 7: USDC
    @{ Row=8; D='0.527'; E='  +0.23%  ' }  # row 8: XRP
    @{ Row=9; D='7.39'; E='  +1.67%  ' }  # row 9: Toncoin
    @{ Row=10; E='  +3.41%  ' }  # row 10: Dogecoin
    @{ Row=11; E='  -1.65%  ' }  # row 11: Cardano
    @{ Row=12; D='3.813.16'; E='  +3.29%  ' }  # row 12: WrappedliquidstakedEther2.0
    @{ Row=13; E='  -2.03%  ' }  # row 13: TRON
    @{ Row=14; D='26.31'; E='  +2.72%  ' }  # row 14: Avalanche
    @{ Row=15; E='  +3.36%  ' }  # row 15: ShibaInu
    @{ Row=16; D='60.942.60'; E='  +4.17%  ' }  # row 16: WrappedBTC
    @{ Row=17; D='3.253.94'; E='  +3.34%  ' }  # row 17: WrappedEther
    @{ Row=18; E='  +1.58%  ' }  # row 18: Polkadot
    @{ Row=19; D='13.43'; E='  +4.00%  ' }  # row 19: Chainlink
    @{ Row=20; D='8.40'; E='  +3.64%  ' }  # row 20: Uniswap
    @{ Row=21; D='377.12'; E='  +1.38%  ' }  # row 21: BitcoinCash
    @{ Row=22; E='  -0.11%  ' }  # row 22: Dai
    @{ Row=23; D='0.531'; E='  +0.58%  ' }  # row 23: Polygon
    @{ Row=24; D='69.98'; E='  +0.59%  ' }  # row 24: Litecoin
    @{ Row=25; D='0.171'; E='  +1.85%  ' }  # row 25: Kaspa
    @{ Row=26; D='8.62'; E='  +3.19%  ' }  # row 26: InternetComputer(DFINITY)
    @{ Row=27; E='  +0.28%  ' }  # row 27: Binance-PegBSC-USD
    @{ Row=28; D='0.0₃0918'; E='  +7.71%  ' }  # row 28: PEPE
    @{ Row=29; E='  +3.23%  ' }  # row 29: PancakeSwap
    @{ Row=30; E='  +1.25%  ' }  # row 30: EthereumClassic
    @{ Row=31; D='6.20'; E='  +3.81%  ' }  # row 31: RenderToken
    @{ Row=32; D='5.43'; E='  +6.06%  ' }  # row 32: NEARProtocol
    @{ Row=33; E='  +7.95%  ' }  # row 33: Fetch.AI
    @{ Row=34; D='6.65'; E='  +5.57%  ' }  # row 34: Aptos
    @{ Row=35; D='158.99'; E='  +1.51%  ' }  # row 35: Monero
    @{ Row=36; D='1.44'; E='  +8.19%  ' }  # row 36: ImmutableX
    @{ Row=37; D='26.40'; E='  +5.95%  ' }  # row 37: EnergySwap
    @{ Row=38; D='2.811.93'; E='  +4.20%  ' }  # row 38: Maker
    @{ Row=39; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0720'; E='  +4.35%  ' }  # row 39: Hedera
    @{ Row=40; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0315'; E='  +7.91%  ' }  # row 40: VeChain
    @{ Row=41; D='1.73'; E='  +2.84%  ' }  # row 41: Stacks
    @{ Row=42; E='  +0.71%  ' }  # row 42: Filecoin
    @{ Row=43; D='40.06'; E='  +2.67%  ' }  # row 43: OKB
    @{ Row=44; D='0.733'; E='  +1.81%  ' }  # row 44: Mantle
    @{ Row=45; D='3.295.24'; E='  +3.18%  ' }  # row 45: RenzoRestakedETH
    @{ Row=46; E='  +2.88%  ' }  # row 46: Stellar
    @{ Row=47; D='1.01'; E='  +3.20%  ' }  # row 47: ONDO
    @{ Row=48; D='21.41'; E='  +7.04%  ' }  # row 48: InjectiveProtocol
    @{ Row=49; E='  +0.90%  ' }  # row 49: Cosmos
    @{ Row=50; D='0.803'; E='  +7.75%  ' }  # row 50: SuiNetwork
    @{ Row=51; D='276.60'; E='  +7.80%  ' }  # row 51: Bittensor
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        # Price column holds free-form text (e.g. "3.249.11"); force text format
        # so Excel doesn't reinterpret it as a number and drop the formatting.
        $ws.Cells.Item($row, 4).NumberFormat = '@'
        $ws.Cells.Item($row, 4).Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
